# Update odds on "Jogos da Semana" sheet to the refreshed FlashScore values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Atletico-MG x Gremio)
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7

# Row 4 (Paysandu PA x Chapecoense-SC)
$ws.Range("G4").Value = 1.62
$ws.Range("I4").Value = 6.25
$ws.Range("J4").Value = 2.25
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("W4").Value = 6
$ws.Range("AC4").Value = 8
$ws.Range("AE4").Value = 19
$ws.Range("AJ4").Value = 21
$ws.Range("AP4").Value = 23
$ws.Range("AQ4").Value = 29

# Row 5 (Envigado x Dep. Pasto)
$ws.Range("I5").Value = 2
$ws.Range("L5").Value = 2.88
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 2.25
$ws.Range("Q5").Value = 2.88
$ws.Range("R5").Value = 1.4
$ws.Range("S5").Value = 1.62
$ws.Range("T5").Value = 2.2
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("W5").Value = 8.5
$ws.Range("AC5").Value = 5.5
$ws.Range("AH5").Value = 5
$ws.Range("AJ5").Value = 10
$ws.Range("AT5").Value = 2.2
$ws.Range("AU5").Value = 10
$ws.Range("BB5").Value = 301

# Row 12 (TNS x Caernarfon)
$ws.Range("G12").Value = 1.16
$ws.Range("H12").Value = 6.9
$ws.Range("I12").Value = 14
$ws.Range("J12").Value = 1.47
$ws.Range("K12").Value = 3.1
$ws.Range("L12").Value = 9.5
$ws.Range("O12").Value = 1.08
$ws.Range("P12").Value = 6.4
$ws.Range("Q12").Value = 1.28
$ws.Range("R12").Value = 3.4
$ws.Range("S12").Value = 1.18
$ws.Range("T12").Value = 4.3
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.91
$ws.Range("W12").Value = 13
$ws.Range("X12").Value = 8.5
$ws.Range("Y12").Value = 10.5
$ws.Range("AB12").Value = 24
$ws.Range("AD12").Value = 16
$ws.Range("AE12").Value = 25
$ws.Range("AG12").Value = 500
$ws.Range("AH12").Value = 55
$ws.Range("AI12").Value = 150
$ws.Range("AJ12").Value = 45
$ws.Range("AK12").Value = 600
$ws.Range("AL12").Value = 175
$ws.Range("AM12").Value = 100
$ws.Range("AO12").Value = 4.75
$ws.Range("AP12").Value = 12
$ws.Range("AQ12").Value = 9.25
$ws.Range("AR12").Value = 24
$ws.Range("AT12").Value = 4.3
$ws.Range("AU12").Value = 9
$ws.Range("AW12").Value = 14
$ws.Range("AX12").Value = 75
$ws.Range("AY12").Value = 50
$ws.Range("AZ12").Value = 500
$ws.Range("BA12").Value = 400
$ws.Range("BB12").Value = 500
